$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Capital One"
$ws.Range("B5").Value = "AIR-INK: Air-Pollution to ink"
